$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - G2: reorder "Recorded By" list
$ws.Range("G2").Value = "System, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg"

# Row 3 - G3: reorder "Recorded By" list
$ws.Range("G3").Value = "System, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 4 - G4: reorder "Recorded By" list
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 5 - G5: reorder "Recorded By" list
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 6 - G6: reorder "Recorded By" list; L6: Recorded Sessions 19 -> 20
$ws.Range("G6").Value = "Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg"
$ws.Range("L6").Value = 20

# Row 7 - G7: reorder "Recorded By" list; L7: Missing Sessions 3 -> 2
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("L7").Value = 2

# Row 9 - G9: reorder "Recorded By" list; L9: Coverage % 65.5% -> 69.0%
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
# Use formula+paste-values so the percent-like text is stored literally
# (not auto-converted to a formatted number) while keeping the cell's style.
$ws.Range("L9").Formula = '="69.0%"'
$ws.Range("L9").Copy()
$ws.Range("L9").PasteSpecial(-4163)

# Row 10 - L10: Average Attendance % 26.2% -> 25.1%
$ws.Range("L10").Formula = '="25.1%"'
$ws.Range("L10").Copy()
$ws.Range("L10").PasteSpecial(-4163)

# Row 12 - G12: reorder "Recorded By" list
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"

# Row 15 - G15: reorder "Recorded By" list; O15/P15/R15/S15 group stats update
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("O15").Value = 20
$ws.Range("P15").Value = 2
$ws.Range("R15").Formula = '="69.0%"'
$ws.Range("R15").Copy()
$ws.Range("R15").PasteSpecial(-4163)
$ws.Range("S15").Formula = '="25.1%"'
$ws.Range("S15").Copy()
$ws.Range("S15").PasteSpecial(-4163)

# Row 17 - session now recorded: style changes from "Not Recorded" (s=5) to default (s=2),
# and G17/H17/I17 filled in
$ws.Range("A17:I17").Style = $ws.Range("A16:I16").Style
$ws.Range("G17").Value = "mohamed.saleem@med.asu.edu.eg"
$ws.Range("H17").Value = "12/251"
$ws.Range("I17").Value = "Recorded"

# Row 28 - G28: reorder "Recorded By" list
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# Row 30 - G30: reorder "Recorded By" list
$ws.Range("G30").Value = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$excel.CutCopyMode = 0
